$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Activités")
$ws2 = $wb.Worksheets.Item("Données")

# --- Sheet1 ("Activités"): fill in row 26 with the new "Chat" entry ---
# (written in this order so new shared strings land at indices 66/67/68,
# matching "22 mai", "Chat", "Nous pouvons maintenant voir ...")
$ws1.Cells.Item(26, 1).Value = "22 mai"
$ws1.Cells.Item(26, 2).Value = 0.33333333333333331
$ws1.Cells.Item(26, 3).Value = 0.5
$ws1.Cells.Item(26, 4).Formula = "=IF(ISBLANK(C26), NOW(),C26)-IF(ISBLANK(B26),NOW(),B26)"
$ws1.Cells.Item(26, 5).Value = "Chat"
$ws1.Cells.Item(26, 6).Value = "Réalisation"
$ws1.Cells.Item(26, 7).Value = "Nous pouvons maintenant voir quel utilisateur est connecté."
$ws1.Rows.Item(26).RowHeight = 30

# --- Sheet2 ("Données"): insert "Chat" alphabetically at row 5, ---
# --- pushing existing rows 5-11 down to 6-12 (column A only). ---
$ws2.Rows.Item(5).Insert()
$ws2.Cells.Item(5, 1).Value = "Chat"

# --- Selections to match the saved view state ---
# (sheet2 selection is set first, then sheet1 is reactivated/selected last so
# that it keeps tabSelected="1" like the original workbook)
$null = $ws2.Range("A5").Select()
$ws1.Activate()
$null = $ws1.Range("H21").Select()
